$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price column so numeric-looking values (e.g. "312.04")
# are stored as text, matching the source data which uses inline strings throughout.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.313.27"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "1.666.34"
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "312.04"
$ws.Range("E5").Value = "  +1.34%  "

$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("D7").Value = "0.3960"
$ws.Range("E7").Value = "  +1.56%  "

$ws.Range("D8").Value = "0.3929"
$ws.Range("E8").Value = "  +1.65%  "

$ws.Range("D9").Value = "51.95"
$ws.Range("E9").Value = "  +3.92%  "

$ws.Range("D10").Value = "1.384"
$ws.Range("E10").Value = "  +2.21%  "

$ws.Range("D11").Value = "1.003"
$ws.Range("E11").Value = "  +0.25%  "

$ws.Range("D12").Value = "0.08563"
$ws.Range("E12").Value = "  -1.45%  "

$ws.Range("E13").Value = "  +2.60%  "

$ws.Range("E14").Value = "  +2.53%  "

$ws.Range("D15").Value = "8.001"
$ws.Range("E15").Value = "  +7.20%  "

$ws.Range("E16").Value = "  +2.95%  "

$ws.Range("D17").Value = "1.667.17"
$ws.Range("E17").Value = "  +4.93%  "

$ws.Range("D18").Value = "95.57"
$ws.Range("E18").Value = "  +0.62%  "

$ws.Range("D19").Value = "0.07013"
$ws.Range("E19").Value = "  +1.61%  "

$ws.Range("D20").Value = "20.46"
$ws.Range("E20").Value = "  -0.56%  "

$ws.Range("D21").Value = "6.992"
$ws.Range("E21").Value = "  +1.23%  "

$ws.Range("D22").Value = "1.002"
$ws.Range("E22").Value = "  +0.09%  "

$ws.Range("D23").Value = "13.80"
$ws.Range("E23").Value = "  +1.47%  "

$ws.Range("D24").Value = "24.322.00"
$ws.Range("E24").Value = "  +0.67%  "

$ws.Range("D25").Value = "2.535"
$ws.Range("E25").Value = "  +8.89%  "

$ws.Range("D26").Value = "3.092"
$ws.Range("E26").Value = "  +10.78%  "

$ws.Range("E27").Value = "  +0.19%  "

$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("D29").Value = "142.21"
$ws.Range("E29").Value = "  +1.01%  "

$ws.Range("D30").Value = "5.383"
$ws.Range("E30").Value = "  +0.56%  "

$ws.Range("D31").Value = "8.012"
$ws.Range("E31").Value = "  -5.70%  "

$ws.Range("D32").Value = "2.520"
$ws.Range("E32").Value = "  +4.39%  "

$ws.Range("D33").Value = "1.854.30"
$ws.Range("E33").Value = "  +1.67%  "

$ws.Range("D34").Value = "1.057"
$ws.Range("E34").Value = "  +10.95%  "

$ws.Range("D35").Value = "0.03089"
$ws.Range("E35").Value = "  +6.31%  "

$ws.Range("D36").Value = "0.08265"
$ws.Range("E36").Value = "  +2.64%  "

$ws.Range("D37").Value = "6.851"
$ws.Range("E37").Value = "  -1.76%  "

$ws.Range("D38").Value = "11.15"
$ws.Range("E38").Value = "  +11.16%  "

$ws.Range("D39").Value = "0.2755"
$ws.Range("E39").Value = "  +2.61%  "

$ws.Range("D40").Value = "0.09287"
$ws.Range("E40").Value = "  +0.88%  "

$ws.Range("B41").Value = "Aptos"
$ws.Range("C41").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D41").Value = "13.76"
$ws.Range("E41").Value = "  +5.57%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.7672"
$ws.Range("E42").Value = "  +1.45%  "

$ws.Range("E43").Value = "  -1.56%  "

$ws.Range("D44").Value = "16.60"
$ws.Range("E44").Value = "  +3.51%  "

$ws.Range("D45").Value = "0.7064"
$ws.Range("E45").Value = "  +2.11%  "

$ws.Range("D46").Value = "2.526"
$ws.Range("E46").Value = "  +2.43%  "

$ws.Range("D47").Value = "4.123"
$ws.Range("E47").Value = "  +0.84%  "

$ws.Range("E48").Value = "  +0.16%  "

$ws.Range("D49").Value = "0.08413"

$ws.Range("D50").Value = "136.42"
$ws.Range("E50").Value = "  +2.17%  "

$ws.Range("D51").Value = "1.261"
$ws.Range("E51").Value = "  -0.37%  "

# Restore default formatting on the Price column (remove the temporary text format)
$ws.Range("D2:D51").ClearFormats()
